$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing
# header cell onto the three new header cells before filling them in.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-44)
$ws.Range("AD2:AD44").Value = 92
$ws.Range("AE2:AE44").Value = 70
$ws.Range("AF2:AF44").Value = 0
